# Update the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) for the refreshed data pull.
#
# Note: several "Price" values are plain decimal-looking strings (e.g.
# "578.42"). Excel auto-coerces those to numbers when assigned directly,
# which would change the cell's stored type away from the original text
# representation. Prefixing with a leading apostrophe forces Excel to
# keep (and display) them as text, matching the source data, while the
# apostrophe itself is not stored as part of the cell's value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.325.25'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.588.49'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''578.42'
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('D6').Value = '''190.73'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').Value = '''0.635'
$ws.Range('E7').Value = '  -2.06%  '
$ws.Range('D8').Value = '3.585.54'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '''0.180'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('D11').Value = '''0.664'
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = '''56.26'
$ws.Range('E12').Value = '  -2.41%  '
$ws.Range('D13').Value = '''0.0000304'
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('D14').Value = '''9.65'
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').Value = '4.166.92'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').Value = '''19.94'
$ws.Range('E16').Value = '  +2.88%  '
$ws.Range('D17').Value = '3.586.17'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').Value = '70.191.38'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').Value = '''1.05'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').Value = '''478.10'
$ws.Range('E22').Value = '  -2.85%  '
$ws.Range('D23').Value = '''19.03'
$ws.Range('E23').Value = '  +13.83%  '
$ws.Range('E24').Value = '  -7.42%  '
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('D26').Value = '''93.60'
$ws.Range('E26').Value = '  +3.31%  '
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('D28').Value = '''11.07'
$ws.Range('E28').Value = '  -1.11%  '
$ws.Range('D29').Value = '''9.39'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = '''32.37'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '''7.73'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('E32').Value = '  +3.54%  '
$ws.Range('D33').Value = '''12.27'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('E34').Value = '  +1.72%  '
$ws.Range('D35').Value = '''582.84'
$ws.Range('E35').Value = '  -5.06%  '
$ws.Range('D36').Value = '''39.27'
$ws.Range('E36').Value = '  +3.55%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '0.0₃0802'
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('D39').Value = '''0.399'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').Value = '''3.31'
$ws.Range('E40').Value = '  +21.58%  '
$ws.Range('E41').Value = '  -5.37%  '
$ws.Range('E42').Value = '  -4.70%  '
$ws.Range('D43').Value = '3.242.86'
$ws.Range('E43').Value = '  -3.63%  '
$ws.Range('E44').Value = '  +7.94%  '
$ws.Range('E45').Value = '  +1.35%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = '''9.58'
$ws.Range('E47').Value = '  +4.99%  '
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '''3.15'
$ws.Range('E51').Value = '  -6.02%  '
